# Refresh the cryptocurrency Price (D) / Volume(1h) (E) columns with the
# latest scraped figures (GitHub Actions "Updated cryptos list" run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "75.888.03"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.906.46"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "596.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.195"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.905.54"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.422"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +13.20%  "
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.440.96"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.697.14"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000189"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.24"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.909.01"
$ws.Range("D18").ClearFormats()
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.71"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.20"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.056.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "499.77"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.69"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.12"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.11"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.24%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "179.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.99"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0911"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.16%  "
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("E47").Value = "  -3.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.575"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.662"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.53%  "
$ws.Range("E51").Value = "  -1.38%  "
